$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-11 from 46079 to 46081
$ws.Range("C2:C11").Value = 46081

# Rows 8-11 got re-sorted: the old row 8 record ("A 25610-2024") moved down to
# row 11, while rows 9, 10, 11 each shifted up by one row. Capture the original
# values first, then write them to their new positions.
$a8 = $ws.Range("A8").Value2
$b8 = $ws.Range("B8").Value2
$g8 = $ws.Range("G8").Value2

$a9 = $ws.Range("A9").Value2
$b9 = $ws.Range("B9").Value2
$g9 = $ws.Range("G9").Value2

$a10 = $ws.Range("A10").Value2
$b10 = $ws.Range("B10").Value2
$g10 = $ws.Range("G10").Value2

$a11 = $ws.Range("A11").Value2
$b11 = $ws.Range("B11").Value2
$g11 = $ws.Range("G11").Value2

$ws.Range("A8").Value = $a9
$ws.Range("B8").Value = $b9
$ws.Range("G8").Value = $g9

$ws.Range("A9").Value = $a10
$ws.Range("B9").Value = $b10
$ws.Range("G9").Value = $g10

$ws.Range("A10").Value = $a11
$ws.Range("B10").Value = $b11
$ws.Range("G10").Value = $g11

$ws.Range("A11").Value = $a8
$ws.Range("B11").Value = $b8
$ws.Range("G11").Value = $g8
